$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking price strings that must remain text (matching the
# original inline-string cell type). Force text format before assigning, then reset
# the cell style back to Normal so no stray style index is left on the cell.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '26.020.56'
$ws.Range("E2").Value = '  -0.72%  '
Set-TextValue $ws.Range("D3") '1.633.74'
$ws.Range("E3").Value = '  -2.29%  '
$ws.Range("E4").Value = '  -0.21%  '
Set-TextValue $ws.Range("D5") '210.89'
$ws.Range("E5").Value = '  -0.33%  '
Set-TextValue $ws.Range("D6") '0.5227'
$ws.Range("E6").Value = '  -1.00%  '
$ws.Range("E7").Value = '  -0.17%  '
Set-TextValue $ws.Range("D8") '0.2585'
$ws.Range("E8").Value = '  -2.12%  '
Set-TextValue $ws.Range("D9") '0.06271'
$ws.Range("E9").Value = '  -0.24%  '
Set-TextValue $ws.Range("D10") '20.55'
$ws.Range("E10").Value = '  -3.59%  '
Set-TextValue $ws.Range("D11") '0.07593'
$ws.Range("E11").Value = '  +0.50%  '
Set-TextValue $ws.Range("D12") '1.626.59'
$ws.Range("E12").Value = '  -2.68%  '
Set-TextValue $ws.Range("D13") '4.423'
$ws.Range("E13").Value = '  -0.55%  '
Set-TextValue $ws.Range("D14") '1.854.52'
$ws.Range("E14").Value = '  -2.43%  '
Set-TextValue $ws.Range("D15") '0.5489'
$ws.Range("E15").Value = '  -2.01%  '
Set-TextValue $ws.Range("D16") '0.0₅8003'
$ws.Range("E16").Value = '  +0.01%  '
Set-TextValue $ws.Range("D17") '64.70'
$ws.Range("E17").Value = '  -3.62%  '
Set-TextValue $ws.Range("D18") '26.007.13'
$ws.Range("E18").Value = '  -0.92%  '
$ws.Range("E19").Value = '  -0.10%  '
Set-TextValue $ws.Range("D20") '4.666'
$ws.Range("E20").Value = '  -2.56%  '
Set-TextValue $ws.Range("D21") '185.41'
$ws.Range("E21").Value = '  -1.09%  '
Set-TextValue $ws.Range("D22") '10.14'
$ws.Range("E22").Value = '  -2.51%  '
Set-TextValue $ws.Range("D23") '6.110'
$ws.Range("E23").Value = '  -1.40%  '
$ws.Range("E24").Value = '  -0.20%  '
Set-TextValue $ws.Range("D25") '145.34'
$ws.Range("E25").Value = '  -2.87%  '
Set-TextValue $ws.Range("D26") '0.1211'
$ws.Range("E26").Value = '  -3.77%  '
Set-TextValue $ws.Range("D27") '7.397'
$ws.Range("E27").Value = '  -2.40%  '
Set-TextValue $ws.Range("D28") '15.64'
$ws.Range("E28").Value = '  -1.99%  '
Set-TextValue $ws.Range("D29") '1.373'
$ws.Range("E29").Value = '  +0.55%  '
Set-TextValue $ws.Range("D30") '0.05914'
$ws.Range("E30").Value = '  -4.30%  '
Set-TextValue $ws.Range("D31") '1.243'
$ws.Range("E31").Value = '  -3.16%  '
Set-TextValue $ws.Range("D32") '3.427'
$ws.Range("E32").Value = '  -2.09%  '
$ws.Range("E33").Value = '  -1.05%  '
Set-TextValue $ws.Range("D34") '1.623'
$ws.Range("E34").Value = '  -0.40%  '
Set-TextValue $ws.Range("D35") '0.9803'
$ws.Range("E35").Value = '  -1.98%  '
$ws.Range("E36").Value = '  -1.07%  '
Set-TextValue $ws.Range("D37") '2.748'
$ws.Range("E37").Value = '  +0.39%  '
$ws.Range("E38").Value = '  -4.76%  '
Set-TextValue $ws.Range("D39") '0.01600'
$ws.Range("E39").Value = '  -1.30%  '
Set-TextValue $ws.Range("D40") '0.8492'
$ws.Range("E40").Value = '  -3.30%  '
$ws.Range("E41").Value = '  -0.16%  '
Set-TextValue $ws.Range("D42") '1.037.89'
$ws.Range("E42").Value = '  -5.53%  '
Set-TextValue $ws.Range("D43") '5.669'
$ws.Range("E43").Value = '  -7.45%  '
Set-TextValue $ws.Range("D44") '100.17'
$ws.Range("E44").Value = '  +0.36%  '
Set-TextValue $ws.Range("D45") '1.782.92'
$ws.Range("E45").Value = '  -2.17%  '
$ws.Range("E46").Value = '  -2.30%  '
Set-TextValue $ws.Range("D47") '54.96'
$ws.Range("E47").Value = '  -1.71%  '
Set-TextValue $ws.Range("D48") '0.9958'
$ws.Range("E48").Value = '  -0.73%  '
Set-TextValue $ws.Range("D49") '8.049'
$ws.Range("E49").Value = '  +0.48%  '
$ws.Range("E50").Value = '  -1.28%  '
Set-TextValue $ws.Range("D51") '0.4225'
$ws.Range("E51").Value = '  -0.68%  '
